$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.537.69'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '2.636.12'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '111.83'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').Value = '326.06'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '0.525'
$ws.Range('E7').Value = '  -0.82%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.548'
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('D10').Value = '39.56'
$ws.Range('E10').Value = '  -3.13%  '
$ws.Range('D11').Value = '20.06'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('E13').Value = '  +1.43%  '
$ws.Range('D14').Value = '7.52'
$ws.Range('E14').Value = '  +2.67%  '
$ws.Range('D15').Value = '3.053.08'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').Value = '2.640.84'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('D18').Value = '49.506.43'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('D19').Value = '13.24'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('D21').Value = '6.67'
$ws.Range('E21').Value = '  -1.43%  '
$ws.Range('D22').Value = '0.0₃0947'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').Value = '268.41'
$ws.Range('E23').Value = '  -3.27%  '
$ws.Range('D24').Value = '69.07'
$ws.Range('E24').Value = '  -4.12%  '
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').Value = '26.05'
$ws.Range('E26').Value = '  -2.42%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  +2.37%  '
$ws.Range('E29').Value = '  -1.42%  '
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('D31').Value = '34.56'
$ws.Range('E31').Value = '  -3.61%  '
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('E33').Value = '  +1.22%  '
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').Value = '19.01'
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('E37').Value = '  +2.97%  '
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('D40').Value = '129.18'
$ws.Range('E40').Value = '  +2.88%  '
$ws.Range('D41').Value = '22.66'
$ws.Range('E41').Value = '  +3.36%  '
$ws.Range('E42').Value = '  -0.92%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = '0.0333'
$ws.Range('E44').Value = '  +5.94%  '
$ws.Range('D45').Value = '2.061.76'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('D47').Value = '2.15'
$ws.Range('E47').Value = '  +8.91%  '
$ws.Range('E48').Value = '  -6.10%  '
$ws.Range('D49').Value = '8.86'
$ws.Range('E49').Value = '  -2.42%  '
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('D51').Value = '58.61'
$ws.Range('E51').Value = '  -1.06%  '
